$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"  (columns: A=(blank), B=Hora_Scrap, C=Hora_Llegada,
#                      D=Línea, E=Minutos, F=Parada, G=Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 08:30:37"
$ws1.Cells.Item(3, 1).Value = "Total filas: 53"

$rows1 = @(
    @("08:30:26", "08:34", "16_SANTA ANA",         4,  "LP1912"),
    @("08:30:26", "08:39", "23_HERNANDEZ",         9,  "LP1912"),
    @("08:30:26", "08:42", "81_EL PELIGRO",        12, "LP1912"),
    @("08:30:26", "08:43", "14_ABASTO",            13, "LP1912"),
    @("08:30:26", "08:48", "16_SANTA ANA",         18, "LP1912"),
    @("08:30:26", "08:53", "10_OLMOS",             23, "LP1912"),
    @("08:30:26", "09:01", "215A_EL PATO",         31, "LP1912"),
    @("08:30:26", "09:03", "11_ETCHEVERRY",        33, "LP1912"),
    @("08:30:26", "09:04", "23_HERNANDEZ",         34, "LP1912"),
    @("08:30:26", "09:10", "16_P MOR-SANTA ANA",   40, "LP1912"),
    @("08:30:26", "09:13", "10_OLMOS",             43, "LP1912"),
    @("08:30:26", "09:17", "27_EL RETIRO",         47, "LP1912"),
    @("08:30:26", "09:21", "26_HERNANDEZ",         51, "LP1912"),
    @("08:30:26", "09:22", "16_SANTA ANA",         52, "LP1912"),
    @("08:30:26", "09:23", "11_ETCHEVERRY",        53, "LP1912"),
    @("08:30:26", "09:32", "15_ABASTO",            62, "LP1912"),
    @("08:30:26", "09:33", "10_OLMOS",             63, "LP1912"),
    @("08:30:26", "09:42", "215C_EL PATO",         72, "LP1912"),
    @("08:30:26", "09:43", "14_ABASTO",            73, "LP1912"),
    @("08:30:26", "09:52", "15_ABASTO",            82, "LP1912"),
    @("08:30:26", "10:03", "11_ETCHEVERRY",        93, "LP1912"),
    @("08:30:26", "10:06", "23_HERNANDEZ",         96, "LP1912")
)

$r = 33
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = "30/12/2025"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215"  (columns: A=(blank), B=Fecha, C=Hora_Scrap,
#                          D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 08:30:37"
$ws2.Cells.Item(3, 1).Value = "Total filas: 8"

$rows2 = @(
    @("30/12/2025", "08:30:26", "09:01", "215A_EL PATO", 31, "LP1912"),
    @("30/12/2025", "08:30:26", "09:42", "215C_EL PATO", 72, "LP1912")
)

$r = 8
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173"  (columns: A=(blank), B=Fecha, C=Hora_Scrap,
#                         D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 08:30:37"
$ws3.Cells.Item(3, 1).Value = "Total filas: 8"

$rows3 = @(
    @("30/12/2025", "08:30:37", "08:38", "215A_LA PLATA",            8,  "L6173"),
    @("30/12/2025", "08:30:32", "09:09", "215D_LA PLATA",            39, "L6203"),
    @("30/12/2025", "08:30:37", "10:03", "215B_LP-P MOR-40 Y 115",   93, "L6173")
)

$r = 7
foreach ($row in $rows3) {
    $ws3.Cells.Item($r, 2).Value = $row[0]
    $ws3.Cells.Item($r, 3).Value = $row[1]
    $ws3.Cells.Item($r, 4).Value = $row[2]
    $ws3.Cells.Item($r, 5).Value = $row[3]
    $ws3.Cells.Item($r, 6).Value = $row[4]
    $ws3.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
